# added some error catching mechanisms
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summoner Name / Summoner Role text updates (columns G/H) ---
$ws.Range("G2").Value = "Halter Penguen"
$ws.Range("H2").Value = "SOLO"

$ws.Range("G3").Value = "Mrs Máster Yi"
$ws.Range("H3").Value = "SOLO"

$ws.Range("G4").Value = "Negabrione"
$ws.Range("H4").Value = "SOLO"

$ws.Range("G5").Value = "HornyCorn"
$ws.Range("H5").Value = "DUO"

$ws.Range("G6").Value = "MagusApex"
$ws.Range("H6").Value = "DUO"

# --- Numeric stat updates (columns B-F, I-K) ---
$ws.Range("B2").Value = 3.07100063734863
$ws.Range("C2").Value = 4818.4
$ws.Range("D2").Value = 0.02179732313575526
$ws.Range("E2").Value = 34.2
$ws.Range("F2").Value = 198
$ws.Range("I2").Value = 0.1261950286806883
$ws.Range("J2").Value = 4.8
$ws.Range("K2").Value = 0.003059273422562141

$ws.Range("B3").Value = 3.833485927892299
$ws.Range("C3").Value = 5897.2
$ws.Range("D3").Value = 0.02445485819541133
$ws.Range("E3").Value = 45.2
$ws.Range("F3").Value = 142
$ws.Range("I3").Value = 0.09318200287130592
$ws.Range("J3").Value = 8.6
$ws.Range("K3").Value = 0.005347130048487147

$ws.Range("B4").Value = 7.829359333360753
$ws.Range("C4").Value = 16197.8
$ws.Range("D4").Value = 0.08013202994381621
$ws.Range("E4").Value = 165
$ws.Range("F4").Value = 430.2
$ws.Range("I4").Value = 0.2063731471137939
$ws.Range("J4").Value = 19.6
$ws.Range("K4").Value = 0.009152606233594018

$ws.Range("B5").Value = 5.64523200920185
$ws.Range("C5").Value = 8954.200000000001
$ws.Range("D5").Value = 0.02991473419815112
$ws.Range("E5").Value = 49.4
$ws.Range("F5").Value = 610.2
$ws.Range("I5").Value = 0.3863354585100985
$ws.Range("J5").Value = 15.2
$ws.Range("K5").Value = 0.009556019006126104

$ws.Range("B6").Value = 3.652230980643308
$ws.Range("C6").Value = 5297.6
$ws.Range("D6").Value = 0.03020442432290065
$ws.Range("E6").Value = 45.8
$ws.Range("F6").Value = 292
$ws.Range("I6").Value = 0.1968040937674682
$ws.Range("J6").Value = 5.4
$ws.Range("K6").Value = 0.00344223402841273
